$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated benchmark data (FreeMarker optimization results) ---
$B = @{6=31977;7=16956;8=11649;9=9833;10=8732;11=7754;12=7269;13=6988;14=6887;15=7096;16=7196;17=6719}
$C = @{6=15639;7=8266;8=6411;9=6182;10=5752;11=5441;12=5422;13=5785;14=5501;15=5667;16=5852;17=5699}

foreach ($r in 6..17) {
    $ws.Cells.Item($r, 2).Value = $B[$r]
    $ws.Cells.Item($r, 3).Value = $C[$r]
}

# --- New column D: percentage improvement of JSLT vs FreeMarker ---
foreach ($r in 6..17) {
    $ws.Cells.Item($r, 4).Formula = "=(B$r-C$r)/B$r"
    $ws.Cells.Item($r, 4).NumberFormat = "0.00%"
}

# --- A4 no longer holds a (blank) cell ---
$ws.Range("A4").ClearContents()

# --- New performance-test summary formulas ---
$ws.Range("B21").Formula = "=1000000/B17*1000"
$ws.Range("C21").Formula = "=1000000/C12*1000"

# --- Reposition the chart to make room for the new column ---
$cht = $ws.ChartObjects().Item(1)
$cht.Left = 294.4875
$cht.Top = 37.35
$cht.Width = 444.2625
$cht.Height = 246.6

# --- Selection moves to C21 ---
$ws.Range("C21").Select()
